$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '29.573.74'
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.66%  '

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.855.03'
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '

# Row 4
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '243.92'
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.51%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '0.6433'
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '1.0000'
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '

# Row 8
$c = $ws.Cells.Item(8, 2)
$c.NumberFormat = "@"
$c.Value = 'Cardano'
$c = $ws.Cells.Item(8, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.3011'
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.53%  '

# Row 9
$c = $ws.Cells.Item(9, 2)
$c.NumberFormat = "@"
$c.Value = 'Dogecoin'
$c = $ws.Cells.Item(9, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.07531'
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.46%  '

# Row 10
$c = $ws.Cells.Item(10, 2)
$c.NumberFormat = "@"
$c.Value = 'Solana'
$c = $ws.Cells.Item(10, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '24.38'
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.74%  '

# Row 11
$c = $ws.Cells.Item(11, 2)
$c.NumberFormat = "@"
$c.Value = 'TRON'
$c = $ws.Cells.Item(11, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.07670'
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.12%  '

# Row 12
$c = $ws.Cells.Item(12, 2)
$c.NumberFormat = "@"
$c.Value = 'WrappedEther'
$c = $ws.Cells.Item(12, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '1.919.34'
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.28%  '

# Row 13
$c = $ws.Cells.Item(13, 2)
$c.NumberFormat = "@"
$c.Value = 'Polkadot'
$c = $ws.Cells.Item(13, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '5.057'
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '

# Row 14
$c = $ws.Cells.Item(14, 2)
$c.NumberFormat = "@"
$c.Value = 'Polygon'
$c = $ws.Cells.Item(14, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.6906'
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.73%  '

# Row 15
$c = $ws.Cells.Item(15, 2)
$c.NumberFormat = "@"
$c.Value = 'Litecoin'
$c = $ws.Cells.Item(15, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '84.06'
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.15%  '

# Row 16
$c = $ws.Cells.Item(16, 2)
$c.NumberFormat = "@"
$c.Value = 'ShibaInu'
$c = $ws.Cells.Item(16, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '0.000009600'
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.60%  '

# Row 17
$c = $ws.Cells.Item(17, 2)
$c.NumberFormat = "@"
$c.Value = 'Uniswap'
$c = $ws.Cells.Item(17, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '6.261'
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.89%  '

# Row 18
$c = $ws.Cells.Item(18, 2)
$c.NumberFormat = "@"
$c.Value = 'WrappedliquidstakedEther2.0'
$c = $ws.Cells.Item(18, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '2.164.76'
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.95%  '

# Row 19
$c = $ws.Cells.Item(19, 2)
$c.NumberFormat = "@"
$c.Value = 'WrappedBTC'
$c = $ws.Cells.Item(19, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '29.625.02'
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '

# Row 20
$c = $ws.Cells.Item(20, 2)
$c.NumberFormat = "@"
$c.Value = 'BitcoinCash'
$c = $ws.Cells.Item(20, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '238.47'
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.60%  '

# Row 21
$c = $ws.Cells.Item(21, 2)
$c.NumberFormat = "@"
$c.Value = 'Avalanche'
$c = $ws.Cells.Item(21, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '12.63'
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '

# Row 22
$c = $ws.Cells.Item(22, 2)
$c.NumberFormat = "@"
$c.Value = 'Dai'
$c = $ws.Cells.Item(22, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '1.000'
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '

# Row 23
$c = $ws.Cells.Item(23, 2)
$c.NumberFormat = "@"
$c.Value = 'Chainlink'
$c = $ws.Cells.Item(23, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '7.722'
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.66%  '

# Row 24
$c = $ws.Cells.Item(24, 2)
$c.NumberFormat = "@"
$c.Value = 'BinanceUSD'
$c = $ws.Cells.Item(24, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '1.000'
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '

# Row 25
$c = $ws.Cells.Item(25, 2)
$c.NumberFormat = "@"
$c.Value = 'Monero'
$c = $ws.Cells.Item(25, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '157.30'
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.02%  '

# Row 26
$c = $ws.Cells.Item(26, 2)
$c.NumberFormat = "@"
$c.Value = 'Stellar'
$c = $ws.Cells.Item(26, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '0.1417'
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.09%  '

# Row 27
$c = $ws.Cells.Item(27, 2)
$c.NumberFormat = "@"
$c.Value = 'Cosmos'
$c = $ws.Cells.Item(27, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '8.541'
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.25%  '

# Row 28
$c = $ws.Cells.Item(28, 2)
$c.NumberFormat = "@"
$c.Value = 'EthereumClassic'
$c = $ws.Cells.Item(28, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '17.85'
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.57%  '

# Row 29
$c = $ws.Cells.Item(29, 2)
$c.NumberFormat = "@"
$c.Value = 'PancakeSwap'
$c = $ws.Cells.Item(29, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '1.489'
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.82%  '

# Row 30
$c = $ws.Cells.Item(30, 2)
$c.NumberFormat = "@"
$c.Value = 'Hedera'
$c = $ws.Cells.Item(30, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '0.05979'
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.98%  '

# Row 31
$c = $ws.Cells.Item(31, 2)
$c.NumberFormat = "@"
$c.Value = 'Toncoin'
$c = $ws.Cells.Item(31, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '1.260'
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.25%  '

# Row 32
$c = $ws.Cells.Item(32, 2)
$c.NumberFormat = "@"
$c.Value = 'Filecoin'
$c = $ws.Cells.Item(32, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '4.153'
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.24%  '

# Row 33
$c = $ws.Cells.Item(33, 2)
$c.NumberFormat = "@"
$c.Value = 'InternetComputer(DFINITY)'
$c = $ws.Cells.Item(33, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '4.080'
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.71%  '

# Row 34
$c = $ws.Cells.Item(34, 2)
$c.NumberFormat = "@"
$c.Value = 'LidoDAOToken'
$c = $ws.Cells.Item(34, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '1.889'
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.44%  '

# Row 35
$c = $ws.Cells.Item(35, 2)
$c.NumberFormat = "@"
$c.Value = 'ARBITRUM'
$c = $ws.Cells.Item(35, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '1.175'
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.37%  '

# Row 36
$c = $ws.Cells.Item(36, 2)
$c.NumberFormat = "@"
$c.Value = 'ImmutableX'
$c = $ws.Cells.Item(36, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.7244'
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.07%  '

# Row 37
$c = $ws.Cells.Item(37, 2)
$c.NumberFormat = "@"
$c.Value = 'HuobiToken'
$c = $ws.Cells.Item(37, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '2.603'
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.18%  '

# Row 38
$c = $ws.Cells.Item(38, 2)
$c.NumberFormat = "@"
$c.Value = 'MXToken'
$c = $ws.Cells.Item(38, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '2.784'
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.52%  '

# Row 39
$c = $ws.Cells.Item(39, 2)
$c.NumberFormat = "@"
$c.Value = 'VeChain'
$c = $ws.Cells.Item(39, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '0.01781'
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.69%  '

# Row 40
$c = $ws.Cells.Item(40, 2)
$c.NumberFormat = "@"
$c.Value = 'Maker'
$c = $ws.Cells.Item(40, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '1.212.93'
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.09%  '

# Row 41
$c = $ws.Cells.Item(41, 2)
$c.NumberFormat = "@"
$c.Value = 'TrustWalletToken'
$c = $ws.Cells.Item(41, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.9172'
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.60%  '

# Row 42
$c = $ws.Cells.Item(42, 2)
$c.NumberFormat = "@"
$c.Value = 'FraxShare'
$c = $ws.Cells.Item(42, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '6.191'
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '

# Row 43
$c = $ws.Cells.Item(43, 2)
$c.NumberFormat = "@"
$c.Value = 'RocketPoolETH'
$c = $ws.Cells.Item(43, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '2.071.12'
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.80%  '

# Row 44
$c = $ws.Cells.Item(44, 2)
$c.NumberFormat = "@"
$c.Value = 'PaxDollar'
$c = $ws.Cells.Item(44, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '

# Row 45
$c = $ws.Cells.Item(45, 2)
$c.NumberFormat = "@"
$c.Value = 'Quant'
$c = $ws.Cells.Item(45, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '102.07'
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '

# Row 46
$c = $ws.Cells.Item(46, 2)
$c.NumberFormat = "@"
$c.Value = 'Aave'
$c = $ws.Cells.Item(46, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '67.54'
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.97%  '

# Row 47
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = '  +4.57%  '

# Row 48
$c = $ws.Cells.Item(48, 2)
$c.NumberFormat = "@"
$c.Value = 'Aptos'
$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '7.386'
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = '  +10.26%  '

# Row 49
$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = "@"
$c.Value = 'TheSandbox'
$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.4068'
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.34%  '

# Row 50
$c = $ws.Cells.Item(50, 2)
$c.NumberFormat = "@"
$c.Value = 'EnergySwap'
$c = $ws.Cells.Item(50, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '9.193'
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.13%  '

# Row 51
$c = $ws.Cells.Item(51, 2)
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c = $ws.Cells.Item(51, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '1.665'
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.09%  '
